# Weekly fruit/vegetable price update:
# Insert a new observation row at row 113 (pushing the existing rows 113-148
# down to 114-149) and populate it with the new week's data for the
# "Terminal La Palmera de La Serena - Berenjena" series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 113..148 down to 114..149, leaving a blank row 113.
$ws.Rows.Item(113).Insert()

# Populate the new row 113 with this week's record.
$ws.Range("A113").Value = 8
$ws.Range("B113").Value = "Terminal La Palmera de La Serena"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44754
$ws.Range("E113").Value = 4
$ws.Range("F113").Value = 100112001
$ws.Range("G113").Value = "Berenjena"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 440
$ws.Range("K113").Value = 9500
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 9750
$ws.Range("N113").Value = "`$/caja 50 unidades"
$ws.Range("O113").Value = "Región de Arica y Parinacota"
$ws.Range("P113").Value = 195
$ws.Range("Q113").Value = 50
$ws.Range("R113").Value = "Hortaliza"
